# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Fri Oct  6 19:51:57 UTC 2023 with GitHub Actions".
#
# Column D ("Price") holds numeric-looking text such as "28.017.81" or
# "0.0484" that uses a dot both as a thousands separator and a decimal
# point, so it must stay plain text -- if Excel is allowed to infer a
# number type it mangles the string (drops trailing zeros, collapses the
# thousands dots). We force text via NumberFormat "@" right before the
# assignment and then clear the format again so the cell keeps the exact
# style it had before (no numFmt override left behind).
#
# Column E ("Volume(1h)") is a plain padded percentage string and Excel
# never tries to reinterpret it, so it can be assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates -- must be written as literal text.
$priceUpdates = @{
    'D2' = '27.936.35'
    'D3' = '1.648.30'
    'D5' = '213.60'
    'D6' = '0.526'
    'D8' = '23.47'
    'D12' = '1.886.12'
    'D13' = '1.654.68'
    'D14' = '4.08'
    'D15' = '0.563'
    'D16' = '65.71'
    'D17' = '27.957.95'
    'D18' = '232.65'
    'D19' = '7.68'
    'D20' = '0.0₃0723'
    'D22' = '10.67'
    'D23' = '4.39'
    'D25' = '152.56'
    'D26' = '6.90'
    'D31' = '0.0483'
    'D32' = '3.36'
    'D33' = '1.453.79'
    'D34' = '3.10'
    'D37' = '0.890'
    'D38' = '0.564'
    'D40' = '0.921'
    'D41' = '69.36'
    'D45' = '2.23'
    'D46' = '1.79'
    'D48' = '1.792.87'
    'D49' = '88.98'
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
    $cell.ClearFormats()
}

# Column E (Volume 1h %) updates -- plain text, no coercion risk.
$volumeUpdates = @{
    'E2' = '  +1.59%  '
    'E3' = '  +1.85%  '
    'E4' = '  -0.05%  '
    'E5' = '  +1.08%  '
    'E6' = '  +0.57%  '
    'E7' = '  -0.06%  '
    'E8' = '  +3.07%  '
    'E9' = '  +1.33%  '
    'E10' = '  +0.28%  '
    'E11' = '  -1.43%  '
    'E12' = '  +2.13%  '
    'E13' = '  +2.04%  '
    'E14' = '  +1.41%  '
    'E15' = '  +2.66%  '
    'E16' = '  +1.14%  '
    'E17' = '  +1.72%  '
    'E18' = '  +1.15%  '
    'E19' = '  +2.03%  '
    'E20' = '  +0.51%  '
    'E21' = '  -0.10%  '
    'E22' = '  +4.64%  '
    'E23' = '  +2.57%  '
    'E24' = '  +4.34%  '
    'E25' = '  +1.99%  '
    'E26' = '  +1.28%  '
    'E28' = '  +0.36%  '
    'E29' = '  -0.01%  '
    'E30' = '  +1.73%  '
    'E31' = '  +0.27%  '
    'E32' = '  +2.90%  '
    'E33' = '  +0.60%  '
    'E34' = '  +1.21%  '
    'E35' = '  +2.04%  '
    'E36' = '  -0.47%  '
    'E37' = '  +3.24%  '
    'E38' = '  +0.60%  '
    'E39' = '  +0.91%  '
    'E40' = '  -1.67%  '
    'E41' = '  +2.19%  '
    'E42' = '  +3.10%  '
    'E43' = '  -0.05%  '
    'E44' = '  +0.37%  '
    'E45' = '  +1.14%  '
    'E46' = '  +5.81%  '
    'E47' = '  -0.90%  '
    'E48' = '  +1.81%  '
    'E49' = '  +3.05%  '
    'E50' = '  +0.45%  '
    'E51' = '  +0.72%  '
}

foreach ($cellRef in $volumeUpdates.Keys) {
    $ws.Range($cellRef).Value = $volumeUpdates[$cellRef]
}

